# Auto-generated edit script: updates Leve market-price / profit
# columns (H-N) across several crafter sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 470.42856
$ws.Range("I9").Value = 510.6
$ws.Range("K9").Value = 510.6
$ws.Range("M9").Value = -341.6

$ws.Range("H28").Value = 1236.75
$ws.Range("I28").Value = 638.1111
$ws.Range("K28").Value = 638.1111
$ws.Range("M28").Value = -153.1111

$ws.Range("H32").Value = 5000
$ws.Range("J32").Value = 5000
$ws.Range("L32").Value = 5000
$ws.Range("N32").Value = -5652

$ws.Range("H41").Value = 752.36365
$ws.Range("J41").Value = 1038.75
$ws.Range("L41").Value = 1038.75
$ws.Range("N41").Value = -1918.75

$ws.Range("H53").Value = 680.1429000000001
$ws.Range("I53").Value = 522.2143
$ws.Range("J53").Value = 996
$ws.Range("K53").Value = 522.2143
$ws.Range("L53").Value = 996
$ws.Range("M53").Value = 114.7857
$ws.Range("N53").Value = -2270

$ws.Range("H88").Value = 1450
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 1450
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H98").Value = 7766520.5
$ws.Range("I98").Value = 8269069.5
$ws.Range("K98").Value = 8269069.5
$ws.Range("M98").Value = -8267571.5

$ws.Range("H107").Value = 16130738
$ws.Range("I107").Value = 10418580
$ws.Range("K107").Value = 10418580
$ws.Range("M107").Value = -10416660

$ws.Range("H122").Value = 7766520.5
$ws.Range("I122").Value = 8269069.5
$ws.Range("K122").Value = 24807208.5
$ws.Range("M122").Value = -24804758.5

$ws.Range("H138").Value = 2843.8352
$ws.Range("I138").Value = 1655.15
$ws.Range("J138").Value = 3178.676
$ws.Range("K138").Value = 4965.450000000001
$ws.Range("L138").Value = 9536.028
$ws.Range("M138").Value = 174.5499999999993
$ws.Range("N138").Value = -19816.028

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7411.3335
$ws.Range("I32").Value = 5417.5537
$ws.Range("K32").Value = 5417.5537
$ws.Range("M32").Value = -5130.5537

$ws.Range("H61").Value = 3226.3958
$ws.Range("I61").Value = 1860.4783
$ws.Range("J61").Value = 4483.04
$ws.Range("K61").Value = 1860.4783
$ws.Range("L61").Value = 4483.04
$ws.Range("M61").Value = -1648.4783
$ws.Range("N61").Value = -4907.04

$ws.Range("H74").Value = 58074.8
$ws.Range("I74").Value = 63025.812
$ws.Range("K74").Value = 63025.812
$ws.Range("M74").Value = -62151.812

$ws.Range("H77").Value = 58074.8
$ws.Range("I77").Value = 63025.812
$ws.Range("K77").Value = 315129.06
$ws.Range("M77").Value = -310761.06

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws.Range("H132").Value = 1398.9111
$ws.Range("I132").Value = 1466.8422
$ws.Range("J132").Value = 1030.1428
$ws.Range("K132").Value = 4400.5266
$ws.Range("L132").Value = 3090.4284
$ws.Range("M132").Value = -1870.5266
$ws.Range("N132").Value = -8150.428400000001

$ws.Range("H136").Value = 3226.3958
$ws.Range("I136").Value = 1860.4783
$ws.Range("J136").Value = 4483.04
$ws.Range("K136").Value = 5581.4349
$ws.Range("L136").Value = 13449.12
$ws.Range("M136").Value = -3031.4349
$ws.Range("N136").Value = -18549.12

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1127.8334
$ws.Range("I22").Value = 553.6
$ws.Range("K22").Value = 553.6
$ws.Range("M22").Value = -380.6

$ws.Range("H99").Value = 2853.1765
$ws.Range("I99").Value = 2853.1765
$ws.Range("K99").Value = 2853.1765
$ws.Range("M99").Value = -1355.1765

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 296126.84
$ws.Range("J31").Value = 3069.5
$ws.Range("L31").Value = 3069.5
$ws.Range("N31").Value = -3659.5

$ws.Range("H34").Value = 296126.84
$ws.Range("J34").Value = 3069.5
$ws.Range("L34").Value = 3069.5
$ws.Range("N34").Value = -3473.5

$ws.Range("H99").Value = 661250.3
$ws.Range("J99").Value = 46666
$ws.Range("L99").Value = 46666
$ws.Range("N99").Value = -49662

$ws.Range("H126").Value = 661250.3
$ws.Range("J126").Value = 46666
$ws.Range("L126").Value = 139998
$ws.Range("N126").Value = -144938

$ws.Range("H132").Value = 4332
$ws.Range("I132").Value = 4005.6428
$ws.Range("K132").Value = 12016.9284
$ws.Range("M132").Value = -9486.928400000001

$ws.Range("H134").Value = 5943.189
$ws.Range("I134").Value = 6390.161
$ws.Range("J134").Value = 3633.8333
$ws.Range("K134").Value = 19170.483
$ws.Range("L134").Value = 10901.4999
$ws.Range("M134").Value = -16635.483
$ws.Range("N134").Value = -15971.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 20502.285
$ws.Range("J36").Value = 22703.2
$ws.Range("L36").Value = 22703.2
$ws.Range("N36").Value = -23673.2

$ws.Range("H132").Value = 346503.66
$ws.Range("I132").Value = 346503.66
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1039510.98
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -1036980.98

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1256.6666
$ws.Range("I22").Value = 872.2222
$ws.Range("K22").Value = 872.2222
$ws.Range("M22").Value = -577.2222

$ws.Range("H27").Value = 1256.6666
$ws.Range("I27").Value = 872.2222
$ws.Range("K27").Value = 872.2222
$ws.Range("M27").Value = -765.2222

$ws.Range("H93").Value = 3000.6365
$ws.Range("I93").Value = 2216.1
$ws.Range("J93").Value = 3654.4167
$ws.Range("K93").Value = 2216.1
$ws.Range("L93").Value = 3654.4167
$ws.Range("M93").Value = -968.0999999999999
$ws.Range("N93").Value = -6150.4167

$ws.Range("H136").Value = 1757.4445
$ws.Range("I136").Value = 1788.5
$ws.Range("K136").Value = 5365.5
$ws.Range("M136").Value = -2815.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H37").Value = 74748.75
$ws.Range("I37").Value = 49499
$ws.Range("J37").Value = 99998.5
$ws.Range("K37").Value = 49499
$ws.Range("L37").Value = 99998.5
$ws.Range("M37").Value = -49296
$ws.Range("N37").Value = -100404.5

$ws.Range("H40").Value = 49999
$ws.Range("I40").Value = 49999
$ws.Range("K40").Value = 49999
$ws.Range("M40").Value = -49850

$ws.Range("H42").Value = 49999
$ws.Range("I42").Value = 49999
$ws.Range("K42").Value = 49999
$ws.Range("M42").Value = -49621

$ws.Range("H132").Value = 1391.5834
$ws.Range("I132").Value = 1391.5834
$ws.Range("K132").Value = 4174.7502
$ws.Range("M132").Value = -1644.7502

$ws.Range("H136").Value = 324004.47
$ws.Range("I136").Value = 334797.47
$ws.Range("K136").Value = 1004392.41
$ws.Range("M136").Value = -1001842.41
